$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "res.groups"
$ws.Range("B6").Value = "base.group_no_one"
$ws.Range("C6").Value = "'True"

$ws.Range("A7").Value = "res.groups"
$ws.Range("B7").Value = "sale.group_delivery_invoice_address"
$ws.Range("C7").Value = "'True"
